$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.460.60"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.093.63"
$ws.Range("E3").Value = "  +9.44%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "250.86"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "0.653"
$ws.Range("E6").Value = "  -6.73%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "47.30"
$ws.Range("E8").Value = "  +6.29%  "
$ws.Range("D9").Value = "59.40"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "0.0740"
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "2.399.30"
$ws.Range("E13").Value = "  +9.47%  "
$ws.Range("D14").Value = "14.43"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D16").Value = "2.090.93"
$ws.Range("E16").Value = "  +9.28%  "
$ws.Range("D17").Value = "5.06"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "36.455.73"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "72.31"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "238.66"
$ws.Range("E22").Value = "  -4.51%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -6.47%  "
$ws.Range("D26").Value = "170.01"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "21.10"
$ws.Range("E27").Value = "  +12.98%  "
$ws.Range("D28").Value = "9.05"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").Value = "  -10.09%  "
$ws.Range("D30").Value = "28.00"
$ws.Range("E30").Value = "  +56.65%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "0.0609"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").Value = "0.0930"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").Value = "0.974"
$ws.Range("E35").Value = "  +10.71%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +14.29%  "
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").Value = "  -6.73%  "
$ws.Range("E40").Value = "  -12.35%  "
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").Value = "96.92"
$ws.Range("E43").Value = "  -8.82%  "
$ws.Range("E44").Value = "  -6.71%  "
$ws.Range("D45").Value = "15.92"
$ws.Range("E45").Value = "  -8.77%  "
$ws.Range("D46").Value = "1.325.17"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("D48").Value = "6.91"
$ws.Range("E48").Value = "  +8.26%  "
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "2.274.29"
$ws.Range("E50").Value = "  +8.52%  "
$ws.Range("D51").Value = "2.21"
$ws.Range("E51").Value = "  -7.00%  "
